# Fruta / hortaliza, semanal
# Inserts one new week of "Papa" price data (Femacal de La Calera) at rows 267-268,
# pushing the existing rows 267..307 down to 269..309.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the block (rows 267 and 268),
# shifting all rows from 267 onward down by two.
$ws.Rows("267:268").Insert()

# --- New row 267 ---
$ws.Cells.Item(267, 1).Value  = 3
$ws.Cells.Item(267, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(267, 3).Value  = 'Coquimbo'
$ws.Cells.Item(267, 4).Value  = 44474
$ws.Cells.Item(267, 5).Value  = 5
$ws.Cells.Item(267, 6).Value  = 100114001
$ws.Cells.Item(267, 7).Value  = 'Papa'
$ws.Cells.Item(267, 8).Value  = 'Asterix'
$ws.Cells.Item(267, 9).Value  = '1a (cosecha)'
$ws.Cells.Item(267, 10).Value = 510
$ws.Cells.Item(267, 11).Value = 9000
$ws.Cells.Item(267, 12).Value = 9500
$ws.Cells.Item(267, 13).Value = 9245
$ws.Cells.Item(267, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(267, 15).Value = 'Región del Maule'
$ws.Cells.Item(267, 16).Value = 370
$ws.Cells.Item(267, 17).Value = 25
$ws.Cells.Item(267, 18).Value = 'Hortaliza'

# --- New row 268 ---
$ws.Cells.Item(268, 1).Value  = 3
$ws.Cells.Item(268, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(268, 3).Value  = 'Coquimbo'
$ws.Cells.Item(268, 4).Value  = 44474
$ws.Cells.Item(268, 5).Value  = 5
$ws.Cells.Item(268, 6).Value  = 100114001
$ws.Cells.Item(268, 7).Value  = 'Papa'
$ws.Cells.Item(268, 8).Value  = 'Rosara'
$ws.Cells.Item(268, 9).Value  = '1a (guarda)'
$ws.Cells.Item(268, 10).Value = 260
$ws.Cells.Item(268, 11).Value = 9000
$ws.Cells.Item(268, 12).Value = 9000
$ws.Cells.Item(268, 13).Value = 9000
$ws.Cells.Item(268, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(268, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(268, 16).Value = 360
$ws.Cells.Item(268, 17).Value = 25
$ws.Cells.Item(268, 18).Value = 'Hortaliza'
